$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: write all cell values in an order that reproduces the target sharedStrings sequence ---
$ws.Range("I7").Value = '>'
$ws.Range("T6").Value = '<'
$ws.Range("L8").Value = 'H'
$ws.Range("J8").Value = 'C'
$ws.Range("H8").Value = 'P'
$ws.Range("AE10").Value = 'prev=H'
$ws.Range("AE11").Value = 'current=head.next'
$ws.Range("AE12").Value = 'head=current.next'
$ws.Range("AD27").Value = 'current.next=prev'
$ws.Range("H7").Value = 1
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = '>'
$ws.Range("L7").Value = 3
$ws.Range("M7").Value = '>'
$ws.Range("N7").Value = 4
$ws.Range("O7").Value = '>'
$ws.Range("P7").Value = 5
$ws.Range("W7").Value = 1
$ws.Range("X7").Value = '>'
$ws.Range("Y7").Value = 2
$ws.Range("I8").Value = '>'
$ws.Range("W8").Value = 'P'
$ws.Range("X8").Value = '<'
$ws.Range("Y8").Value = 'C'
$ws.Range("AA8").Value = 'H'
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = '>'
$ws.Range("J13").Value = 2
$ws.Range("K13").Value = '>'
$ws.Range("L13").Value = 3
$ws.Range("M13").Value = '>'
$ws.Range("N13").Value = 4
$ws.Range("O13").Value = '>'
$ws.Range("P13").Value = 5
$ws.Range("H14").Value = 'P'
$ws.Range("I14").Value = '<'
$ws.Range("J14").Value = 'C'
$ws.Range("L14").Value = 'H'
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = '>'
$ws.Range("J16").Value = 2
$ws.Range("K16").Value = '>'
$ws.Range("L16").Value = 3
$ws.Range("M16").Value = '>'
$ws.Range("N16").Value = 4
$ws.Range("O16").Value = '>'
$ws.Range("P16").Value = 5
$ws.Range("H17").Value = 'P'
$ws.Range("I17").Value = '<'
$ws.Range("J17").Value = 'C'
$ws.Range("L17").Value = 'H'
$ws.Range("J18").Value = 'P'
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = '>'
$ws.Range("J20").Value = 2
$ws.Range("K20").Value = '>'
$ws.Range("L20").Value = 3
$ws.Range("M20").Value = '>'
$ws.Range("N20").Value = 4
$ws.Range("O20").Value = '>'
$ws.Range("P20").Value = 5
$ws.Range("I21").Value = '<'
$ws.Range("J21").Value = 'P'
$ws.Range("L21").Value = 'H'
$ws.Range("L22").Value = 'C'
$ws.Range("H24").Value = 1
$ws.Range("I24").Value = '>'
$ws.Range("J24").Value = 2
$ws.Range("K24").Value = '>'
$ws.Range("L24").Value = 3
$ws.Range("M24").Value = '>'
$ws.Range("N24").Value = 4
$ws.Range("O24").Value = '>'
$ws.Range("P24").Value = 5
$ws.Range("I25").Value = '<'
$ws.Range("J25").Value = 'P'
$ws.Range("K25").Value = '>'
$ws.Range("L25").Value = 'C'
$ws.Range("N25").Value = 'H'

# --- Step 2: center alignment (matches col-level default style already, but set explicitly for clarity) ---

# --- Step 3: apply fills. Blue group (s=2) first, then orange group (s=3), to mirror authoring order. ---
$blueCells = @("I7", "K7", "M7", "O7", "X7", "I8", "I13", "K13", "M13", "O13", "I16", "K16", "M16", "O16", "I20", "K20", "M20", "O20", "I24", "K24", "M24", "O24", "K25")
foreach ($addr in $blueCells) {
    $ws.Range($addr).Interior.Color = 0xF3E3DA
    $ws.Range($addr).HorizontalAlignment = -4108
}

$orangeCells = @("T6", "X8", "I14", "I17", "I21", "I25")
foreach ($addr in $orangeCells) {
    $ws.Range($addr).Interior.Color = 0xADCBF8
    $ws.Range($addr).HorizontalAlignment = -4108
}

# --- Step 4: selection ---
$ws.Range("AD11").Select()
